# Add a "Save" column (column H) to the s_vals sheet, matching the
# formatting of the existing header cells (e.g. G1) and filling the
# data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, styled like the other header cells (bold, bordered,
# centered/top-aligned) by copying the format from G1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data cells H2:H3 with value 0.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
